$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param($ws, $r, $timestamp, $value, $isNumeric)
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-01"
    $ws.Cells.Item($r, 2).Value = $timestamp
    $ws.Cells.Item($r, 3).Value = "20:00"
    $ws.Cells.Item($r, 4).Value = "Bedroom"
    if ($isNumeric) {
        $ws.Cells.Item($r, 5).Value = $value
    } else {
        $ws.Cells.Item($r, 5).Value = [string]$value
    }
    $ws.Cells.Item($r, 6).Value = "Occupied"
}

$wsInBed = $wb.Worksheets.Item("mmWave(InBed)")
Add-LogRow $wsInBed 25 "20:13:06" "In Bed" $false
Add-LogRow $wsInBed 26 "20:15:37" "In Bed" $false
Add-LogRow $wsInBed 27 "20:15:38" "In Bed" $false
Add-LogRow $wsInBed 28 "20:15:39" "In Bed" $false
Add-LogRow $wsInBed 29 "20:15:41" "In Bed" $false
Add-LogRow $wsInBed 30 "20:15:42" "In Bed" $false
Add-LogRow $wsInBed 31 "20:15:47" "In Bed" $false
Add-LogRow $wsInBed 32 "20:15:48" "In Bed" $false
Add-LogRow $wsInBed 33 "20:15:50" "In Bed" $false
Add-LogRow $wsInBed 34 "20:15:51" "In Bed" $false
Add-LogRow $wsInBed 35 "20:15:54" "In Bed" $false
Add-LogRow $wsInBed 36 "20:15:58" "In Bed" $false
Add-LogRow $wsInBed 37 "20:16:03" "In Bed" $false
Add-LogRow $wsInBed 38 "20:16:05" "In Bed" $false

$wsBR = $wb.Worksheets.Item("mmWave(BR)")
Add-LogRow $wsBR 23 "20:15:39" 2 $true
Add-LogRow $wsBR 24 "20:15:40" 11 $true
Add-LogRow $wsBR 25 "20:15:41" 15 $true
Add-LogRow $wsBR 26 "20:15:43" 2 $true
Add-LogRow $wsBR 27 "20:15:48" 25 $true
Add-LogRow $wsBR 28 "20:15:49" 2 $true
Add-LogRow $wsBR 29 "20:15:51" 7 $true
Add-LogRow $wsBR 30 "20:15:52" 2 $true
Add-LogRow $wsBR 31 "20:15:55" 1 $true
Add-LogRow $wsBR 32 "20:15:59" 2 $true
Add-LogRow $wsBR 33 "20:16:04" 100 $true
Add-LogRow $wsBR 34 "20:16:06" 2 $true

$wsHR = $wb.Worksheets.Item("mmWave(HR)")
Add-LogRow $wsHR 23 "20:15:38" 50 $true
Add-LogRow $wsHR 24 "20:15:40" 59 $true
Add-LogRow $wsHR 25 "20:15:41" 63 $true
Add-LogRow $wsHR 26 "20:15:42" 50 $true
Add-LogRow $wsHR 27 "20:15:48" 73 $true
Add-LogRow $wsHR 28 "20:15:49" 50 $true
Add-LogRow $wsHR 29 "20:15:50" 55 $true
Add-LogRow $wsHR 30 "20:15:51" 50 $true
Add-LogRow $wsHR 31 "20:15:55" 49 $true
Add-LogRow $wsHR 32 "20:15:58" 50 $true
Add-LogRow $wsHR 33 "20:16:04" 148 $true
Add-LogRow $wsHR 34 "20:16:05" 50 $true

Write-Output "done"